$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.051.34'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '2.057.57'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.671'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.48'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.79%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.38'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.386'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("E11").Value = '  +6.10%  '
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '16.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.27%  '
$ws.Range("D14").Value = '2.361.44'
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.805'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.32%  '
$ws.Range("D17").Value = '2.064.11'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '37.044.34'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +14.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").Value = '0.0₃0918'
$ws.Range("E21").Value = '  +7.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.30%  '
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("E31").Value = '  +4.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0622'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0886'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +17.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("E45").Value = '  +1.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.68%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.292.01'
$ws.Range("E48").Value = '  -2.93%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -22.57%  '
